$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert two new rows at row 4 (existing row 4 "note" and below shift down to 6+)
$ws.Rows.Item(4).Insert() | Out-Null
$ws.Rows.Item(4).Insert() | Out-Null

# Row 4: dataset.preview.table
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

# Row 5: dataset.preview.line
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# Style rows 4:5 - same base look (vertical center) plus wrap text, and taller rows
$r = $ws.Range("A4:B5")
$r.VerticalAlignment = -4108
$r.WrapText = $true
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120

# Update the visible selection to match new layout
$ws.Range("B7").Select() | Out-Null

Write-Output "done"
